$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2170542635658915
$ws.Range("C2").Value = 0.5193798449612403
$ws.Range("J2").Value = 0.02583979328165375
$ws.Range("P2").Value = 0.1679586563307494
$ws.Range("S2").Value = 0.06976744186046512

# Row 3
$ws.Range("B3").Value = 0.009259259259259259
$ws.Range("C3").Value = 0.01851851851851852
$ws.Range("J3").Value = 0.04166666666666666
$ws.Range("O3").Value = 0.009259259259259259
$ws.Range("P3").Value = 0.7407407407407407
$ws.Range("S3").Value = 0.1805555555555556

# Row 4
$ws.Range("J4").Value = 0.09302325581395349
$ws.Range("P4").Value = 0.5348837209302325
$ws.Range("S4").Value = 0.3720930232558139

# Row 6
$ws.Range("B6").Value = 0.03846153846153846
$ws.Range("D6").Value = 0.01923076923076923
$ws.Range("F6").Value = 0.05288461538461538
$ws.Range("J6").Value = 0.3076923076923077
$ws.Range("O6").Value = 0.01923076923076923
$ws.Range("Q6").Value = 0.1634615384615385
$ws.Range("R6").Value = 0.0625
$ws.Range("S6").Value = 0.3365384615384616

# Row 7
$ws.Range("B7").Value = 0.1392405063291139
$ws.Range("D7").Value = 0.006329113924050633
$ws.Range("F7").Value = 0.05696202531645569
$ws.Range("J7").Value = 0.1518987341772152
$ws.Range("O7").Value = 0.0189873417721519
$ws.Range("Q7").Value = 0.1708860759493671
$ws.Range("R7").Value = 0.06329113924050633
$ws.Range("S7").Value = 0.3924050632911392

# Row 8
$ws.Range("B8").Value = 0.1279620853080569
$ws.Range("D8").Value = 0.01658767772511848
$ws.Range("E8").Value = 0.002369668246445498
$ws.Range("F8").Value = 0.05213270142180094
$ws.Range("J8").Value = 0.1042654028436019
$ws.Range("O8").Value = 0.01895734597156398
$ws.Range("Q8").Value = 0.2014218009478673
$ws.Range("R8").Value = 0.0924170616113744
$ws.Range("S8").Value = 0.3838862559241706

# Row 9
$ws.Range("B9").Value = 0.1451612903225807
$ws.Range("D9").Value = 0.02150537634408602
$ws.Range("F9").Value = 0.08064516129032258
$ws.Range("J9").Value = 0.1075268817204301
$ws.Range("O9").Value = 0.02150537634408602
$ws.Range("Q9").Value = 0.1827956989247312
$ws.Range("R9").Value = 0.07526881720430108
$ws.Range("S9").Value = 0.3655913978494624

# Row 10
$ws.Range("B10").Value = 0.1375579598145286
$ws.Range("D10").Value = 0.02163833075734158
$ws.Range("E10").Value = 0.0007727975270479134
$ws.Range("F10").Value = 0.06182380216383308
$ws.Range("J10").Value = 0.1306027820710974
$ws.Range("O10").Value = 0.02241112828438949
$ws.Range("Q10").Value = 0.1970633693972179
$ws.Range("R10").Value = 0.0695517774343122
$ws.Range("S10").Value = 0.3585780525502318

# Row 11
$ws.Range("G11").Value = 0.1586715867158671
$ws.Range("J11").Value = 0.1180811808118081
$ws.Range("K11").Value = 0.2583025830258303
$ws.Range("L11").Value = 0.4538745387453875
$ws.Range("S11").Value = 0.01107011070110701

# Row 12
$ws.Range("G12").Value = 0.6357142857142857
$ws.Range("J12").Value = 0.2
$ws.Range("K12").Value = 0.007142857142857143
$ws.Range("L12").Value = 0.1142857142857143
$ws.Range("S12").Value = 0.04285714285714286

# Row 13
$ws.Range("G13").Value = 0.6136363636363636
$ws.Range("J13").Value = 0.3181818181818182
$ws.Range("S13").Value = 0.06818181818181818

# Row 14
$ws.Range("G14").Value = 0.8
$ws.Range("J14").Value = 0.2

# Row 15
$ws.Range("F15").Value = 0.004672897196261682
$ws.Range("H15").Value = 0.1495327102803738
$ws.Range("I15").Value = 0.0794392523364486
$ws.Range("J15").Value = 0.3738317757009346
$ws.Range("K15").Value = 0.06074766355140187
$ws.Range("M15").Value = 0.004672897196261682
$ws.Range("O15").Value = 0.03271028037383177
$ws.Range("S15").Value = 0.294392523364486

# Row 16
$ws.Range("H16").Value = 0.2083333333333333
$ws.Range("I16").Value = 0.09583333333333334
$ws.Range("J16").Value = 0.4166666666666667
$ws.Range("K16").Value = 0.09166666666666666
$ws.Range("M16").Value = 0.0125
$ws.Range("N16").Value = 0.0125
$ws.Range("O16").Value = 0.05416666666666667
$ws.Range("S16").Value = 0.1083333333333333

# Row 17
$ws.Range("F17").Value = 0.02277904328018223
$ws.Range("H17").Value = 0.1526195899772209
$ws.Range("I17").Value = 0.08200455580865604
$ws.Range("J17").Value = 0.4533029612756264
$ws.Range("K17").Value = 0.0888382687927107
$ws.Range("M17").Value = 0.02277904328018223
$ws.Range("O17").Value = 0.06378132118451026
$ws.Range("S17").Value = 0.1138952164009112

# Row 18
$ws.Range("F18").Value = 0.01176470588235294
$ws.Range("H18").Value = 0.1470588235294118
$ws.Range("I18").Value = 0.09411764705882353
$ws.Range("J18").Value = 0.4411764705882353
$ws.Range("K18").Value = 0.09411764705882353
$ws.Range("M18").Value = 0.01176470588235294
$ws.Range("N18").Value = 0.01176470588235294
$ws.Range("O18").Value = 0.06470588235294118
$ws.Range("S18").Value = 0.1235294117647059

# Row 19
$ws.Range("F19").Value = 0.01954120645709431
$ws.Range("H19").Value = 0.2098555649957519
$ws.Range("I19").Value = 0.08581138487680544
$ws.Range("J19").Value = 0.3789294817332201
$ws.Range("K19").Value = 0.0994052676295667
$ws.Range("M19").Value = 0.02633814783347494
$ws.Range("O19").Value = 0.09770603228547153
$ws.Range("S19").Value = 0.09770603228547153

Write-Output "Applied 112 cell updates"
